$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the width/value of the time A2 cell
$ws.Range("A2").Value = -8.932

# Set the active selection to A2 (matches saved selection in the edit)
$ws.Range("A2").Select()
